$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '43.836.96'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '2.236.43'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.63'
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.16'
$ws.Range("E6").Value = '  -7.61%  '
$ws.Range("E7").Value = '  -3.21%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  -7.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.75'
$ws.Range("E10").Value = '  -8.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0820'
$ws.Range("E11").Value = '  -2.76%  '
$ws.Range("E12").Value = '  -7.64%  '
$ws.Range("E13").Value = '  -2.82%  '
$ws.Range("D14").Value = '2.576.14'
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D15").Value = '2.237.12'
$ws.Range("E15").Value = '  -2.35%  '
$ws.Range("E16").Value = '  -5.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.89'
$ws.Range("E17").Value = '  -5.31%  '
$ws.Range("D18").Value = '43.681.65'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.08'
$ws.Range("E19").Value = '  -7.81%  '
$ws.Range("D20").Value = '0.0₃0966'
$ws.Range("E20").Value = '  -3.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.28'
$ws.Range("E21").Value = '  -4.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.99'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.03'
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.97'
$ws.Range("E24").Value = '  -7.79%  '
$ws.Range("E25").Value = '  -8.89%  '
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.06'
$ws.Range("E27").Value = '  -2.39%  '
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.54'
$ws.Range("E29").Value = '  -6.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.96'
$ws.Range("E30").Value = '  -9.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.00'
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.99'
$ws.Range("E32").Value = '  -5.03%  '
$ws.Range("E33").Value = '  -6.80%  '
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.90'
$ws.Range("E36").Value = '  -8.90%  '
$ws.Range("E37").Value = '  -6.45%  '
$ws.Range("E38").Value = '  -3.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.53'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("E40").Value = '  -11.43%  '
$ws.Range("E41").Value = '  -12.02%  '
$ws.Range("E42").Value = '  -6.51%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '1.699.21'
$ws.Range("E44").Value = '  -4.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '82.32'
$ws.Range("E45").Value = '  -4.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.194'
$ws.Range("E46").Value = '  -7.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.14'
$ws.Range("E47").Value = '  -5.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.53'
$ws.Range("E48").Value = '  -3.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '70.95'
$ws.Range("E49").Value = '  -6.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '56.18'
$ws.Range("E50").Value = '  -6.11%  '
$ws.Range("E51").Value = '  -5.46%  '
